# Se corrigio color del header - Se agrego ejemplo en los archivos excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New example row: phone number + full name
$ws.Range("A2").Value = 932000076
$ws.Range("B2").Value = "Marco Alex Martinez Ramirez"

# Right-align the new name cell
$ws.Range("B2").HorizontalAlignment = -4152  # xlRight

# Column widths to fit the new content (approx. 15.3 / 29.6 characters)
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 28.6667

# Match the saved selection state from the source workbook
[void]$ws.Range("B12").Select()
